$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "67.363.79"
$ws.Range("E2").Value = "  -3.12%  "

$ws.Range("D3").Value = "3.499.79"
$ws.Range("E3").Value = "  -4.58%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.99%  "

$ws.Range("D7").Value = "3.498.96"
$ws.Range("E7").Value = "  -4.52%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.480"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.32%  "

$ws.Range("E10").Value = "  -4.65%  "

$ws.Range("E11").Value = "  -4.06%  "

$ws.Range("E12").Value = "  -4.26%  "

$ws.Range("E13").Value = "  -5.63%  "

$ws.Range("D14").Value = "4.091.54"
$ws.Range("E14").Value = "  -4.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.37%  "

$ws.Range("D16").Value = "3.492.90"
$ws.Range("E16").Value = "  -5.10%  "

$ws.Range("D17").Value = "67.266.36"
$ws.Range("E17").Value = "  -3.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.117"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("E19").Value = "  -1.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "446.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -13.15%  "

$ws.Range("E23").Value = "  -4.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.95%  "

$ws.Range("E25").Value = "  +5.67%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").Value = "3.639.19"
$ws.Range("E27").Value = "  -4.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.17%  "

$ws.Range("E30").Value = "  -5.46%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.38%  "

$ws.Range("E33").Value = "  +1.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.68"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.11%  "

$ws.Range("D36").Value = "3.489.71"
$ws.Range("E36").Value = "  -4.96%  "

$ws.Range("E37").Value = "  -6.49%  "

$ws.Range("E38").Value = "  -3.68%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "174.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.38%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.51%  "

$ws.Range("E43").Value = "  -2.02%  "

$ws.Range("E44").Value = "  -6.93%  "

$ws.Range("E45").Value = "  -4.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.04%  "

$ws.Range("E47").Value = "  +6.51%  "

$ws.Range("E48").Value = "  -7.12%  "

$ws.Range("E49").Value = "  -5.45%  "

$ws.Range("E50").Value = "  -4.28%  "

$ws.Range("E51").Value = "  -3.93%  "
